# Updates the "Comp controls" sheet in the Miami.xlsx workbook:
#  - Pregate / rename marker & FCS filenames to reflect "-A" (area) channel naming
#    and corrected fluorochrome names (e.g. PE Cy7 -> PE Cy7 YG, Alexa 647/APC-H7 -> APC-A/APC-Cy7-A).
#  - Updates column widths / best-fit and the saved selection to match the refreshed sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comp controls")

# --- FCS file name (column A) and Marker (column B) updates for rows 4-19 ---

$updates = @(
    @{ Row = 4;  A = "Compensation Controls_Live Green FITC-A Stained Control_F02.fcs";              B = "Live Green:FITC-A" },
    @{ Row = 5;  A = "Compensation Controls_CD197 PE-A Stained Control_F03.fcs";                      B = "CD197:PE-A" },
    @{ Row = 6;  A = "Compensation Controls_CD4 PerCP-Cy5-5-A Stained Control_F04.fcs";               B = "CD4:PerCP-Cy5-5-A" },
    @{ Row = 7;  A = "Compensation Controls_CD45RA PE Cy7 YG-A Stained Control_F05.fcs";              B = "CD45RA:PE Cy7 YG-A" },
    @{ Row = 8;  A = "Compensation Controls_CD194 PE Cy7 YG-A Stained Control_F06.fcs";               B = "CD194:PE Cy7 YG-A" },
    @{ Row = 9;  A = "Compensation Controls_CD27 PE Cy7 YG-A Stained Control_F07.fcs";                B = "CD27:PE Cy7 YG-A" },
    @{ Row = 10; A = "Compensation Controls_CD11c PE Cy7 YG-A Stained Control_F08.fcs";               B = "CD11c:PE Cy7 YG-A" },
    @{ Row = 11; A = "Compensation Controls_CD196 PE Cy7 YG-A Stained Control_F09.fcs";               B = "CD196:PE Cy7 YG-A" },
    @{ Row = 12; A = "Compensation Controls_CD38 APC-A Stained Control_F10.fcs";                      B = "CD38:APC-A" },
    @{ Row = 13; A = "Compensation Controls_CD127 APC-A Stained Control_F11.fcs";                     B = "CD127:APC-A" },
    @{ Row = 14; A = "Compensation Controls_CD8 APC-Cy7-A Stained Control_F12.fcs";                   B = "CD8:APC-Cy7-A" },
    @{ Row = 15; A = "Compensation Controls_CD45RO APC-Cy7-A Stained Control_G01.fcs";                B = "CD45RO:APC-Cy7-A" },
    @{ Row = 16; A = "Compensation Controls_CD20 APC-Cy7-A Stained Control_G02.fcs";                  B = "CD20:APC-Cy7-A" },
    @{ Row = 17; A = "Compensation Controls_CD3+CD19+CD20+ APC-Cy7-A Stained Control_G03.fcs";        B = "CD3+19+20:APC-Cy7-A" },
    @{ Row = 18; A = "Compensation Controls_CD3 Pacific Blue-A Stained Control_G04.fcs";              B = "CD3:Pacific Blue-A" },
    @{ Row = 19; A = "Compensation Controls_HLA-DR AmCyan-A Stained Control_G05.fcs";                 B = "HLA-DR:AmCyan-A" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 1).Value = $u.A
    $ws.Cells.Item($u.Row, 2).Value = $u.B
}

# --- Column widths / best fit to accommodate the longer text ---
# (Excel quantizes ColumnWidth to whole pixels internally, so these inputs are
# chosen to land as close as possible to the saved widths of 63.6640625 / 19.83203125.)
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(1).ColumnWidth = 62.75
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(2).ColumnWidth = 18.92

# --- Refresh the stored selection on the sheet ---
$ws.Range("B26").Select()

$wb.Save()
